$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 103. This shifts the existing rows 103:227
# down to 104:228 (carrying their formatting/styles with them), matching
# the target diff where every row from 103 onward ends up holding the
# data that used to belong to the row above it, and a brand new row of
# data is introduced at position 103.
$ws.Rows("103:103").Insert()

# Populate the newly inserted row 103 with its data.
$ws.Range("A103").Value2 = 7
$ws.Range("B103").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C103").Value2 = "Ñuble"
$ws.Range("D103").Value2 = 44789
$ws.Range("E103").Value2 = 16
$ws.Range("F103").Value2 = 100112017
$ws.Range("G103").Value2 = "Apio"
$ws.Range("H103").Value2 = "Americana (o)"
$ws.Range("I103").Value2 = "Primera"
$ws.Range("J103").Value2 = 60
$ws.Range("K103").Value2 = 10000
$ws.Range("L103").Value2 = 10000
$ws.Range("M103").Value2 = 10000
$ws.Range("N103").Value2 = "$/docena de matas"
$ws.Range("O103").Value2 = "Provincia del Elquí"
$ws.Range("P103").Value2 = 1667
$ws.Range("Q103").Value2 = 6
$ws.Range("R103").Value2 = "Hortaliza"
